# Weekly refresh of the "Betarraga" (Hortaliza) sheet:
# a new week's record is inserted at row 442, pushing the previously
# recorded rows 442-475 down to 443-476 (the existing data is otherwise
# untouched - columns A/B/C/E/F/G/H/N/Q/R are constant across all of these
# rows for this market/category/quality-grade combination).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data block (rows 442-475) down by one row.
$ws.Rows.Item(442).Insert()

# Populate the newly inserted row with this week's figures.
$ws.Range("A442").Value = 11
$ws.Range("B442").Value = "Vega Monumental Concepción"
$ws.Range("C442").Value = "Bíobío"
$ws.Range("D442").Value = 45013
$ws.Range("E442").Value = 8
$ws.Range("F442").Value = 100114014
$ws.Range("G442").Value = "Betarraga"
$ws.Range("H442").Value = "Sin especificar"
$ws.Range("I442").Value = "Primera"
$ws.Range("J442").Value = 260
$ws.Range("K442").Value = 600
$ws.Range("L442").Value = 650
$ws.Range("M442").Value = 623
$ws.Range("N442").Value = "$/paquete 5 unidades"
$ws.Range("O442").Value = "Región Metropolitana"
$ws.Range("P442").Value = 125
$ws.Range("Q442").Value = 5
$ws.Range("R442").Value = "Hortaliza"
